# "Loan RBI, Variable Instalments"
#
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# immediately before column N. This shifts the former "Late" / heading /
# "Outstanding" columns one slot to the right (N->O, O->P, P->Q) and
# leaves the freshly inserted column N blank. The sheet also becomes the
# active/selected tab, with cell Q5 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N, pushing existing data right.
$ws.Columns("N:N").Insert()

# Keep the new column's width close to its neighbours (column M).
$ws.Columns("N:N").ColumnWidth = 9.85

# Make "Repayment schedule" the active sheet/tab (was "NewLoanInput").
$ws.Activate()

# Restore the selected cell on the now-active sheet.
[void]$ws.Range("Q5").Select()
